$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data for rows 2-13 (Common Word, Total Frequency, Websites)
$data = @(
    @("Cookies",   83,  "https://www.100-dakar.com (14), https://www.benjaminwahl.at (18), https://www.drehorgelkabarett.at (14), https://www.ottosaxinger.at (3), https://www.peligro.at (14), https://www.schuledesungehorsams.at (2), https://www.skodone.at (18)"),
    @("Seite",     35,  "https://www.hungaromedia.at (8), https://www.kuenstlerinnen.at (8), https://www.luckeneder-art.at (8), https://www.platform-socialism.org (3), https://www.regional-express.org (8)"),
    @("Linz",      52,  "https://www.freie-medien.at (10), https://www.freizeitundkommunikation.at (3), https://www.linzfmr.at (18), https://www.steingeschichten.at (21)"),
    @("März",      49,  "https://www.das-kollektiv.at (8), https://www.feminismus-krawall.at (16), https://www.fiftitu.at (19), https://www.unkraut-comics.at (6)"),
    @("Art",       29,  "https://www.eipcp.net (13), https://www.kairus.org (6), https://www.negentropy-sport.net (2), https://www.radical-openness.org (8)"),
    @("Schule",    294, "https://www.alteschule-gutau.at (4), https://www.derschueler.at (5), https://www.die-schule.at (285)"),
    @("Kultur",    36,  "https://www.frauenkultur.at (17), https://www.igkultur.at (14), https://www.interregnum.live (5)"),
    @("Uhr",       35,  "https://www.fro.at (23), https://www.rudolfhabringer.at (12)"),
    @("Film",      15,  "https://www.corpushomini.info (3), https://www.doublehappiness.at (12)"),
    @("Magdalena", 12,  "https://www.magdalenareiter.at (2), https://www.themagdalenaproject.org (10)"),
    @("Arbeit",    9,   "https://www.eliot.at (7), https://www.hoerstadt.at (2)"),
    @("Andreas",   8,   "https://www.andreaskurz.net (2), https://www.andreaszingerle.com (6)")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove the now-unused rows (previously rows 14-17: Arbeit, Kasino, Andreas, Anna)
$ws.Range("A14:A17").EntireRow.Delete()
